$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.822.32'
$ws.Range('E2').Value = '  -0.41%  '
$ws.Range('D3').Value = '1.637.53'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.97'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.63%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06437'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('E10').Value = '  -0.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07775'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.282'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '1.862.81'
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '1.635.32'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5650'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.50%  '
$ws.Range('D16').Value = '0.0₅7601'
$ws.Range('E16').Value = '  -1.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.19'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.52%  '
$ws.Range('D18').Value = '25.857.59'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '195.81'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.332'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.900'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.098'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.791'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1276'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '139.80'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.25%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.800'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.49'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.243'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04884'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.302'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.225'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.64%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('E35').Value = '  -0.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9046'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.578'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('D38').Value = '1.130.25'
$ws.Range('E38').Value = '  +0.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5509'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01564'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9956'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.532'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8010'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.86'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('D45').Value = '1.772.71'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  -6.81%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.43'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.53%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4433'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E49').Value = '  -2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.652'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.002'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.12%  '
